$d = $word.ActiveDocument

# 1) Replace the document title text.
$d.Content.Find.Execute("2.2 - Debate I", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Placeholder - Check Back Later", 2)

# 2) Remove the trailing " :::" runs that follow "...general edification later."
$d.Content.Find.Execute("general edification later. :::", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "general edification later.", 2)
